# Add "time_taken" metadata column (F) to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, matching the style of the existing header row (B1:E1).
$ws.Cells.Item(1, 6).Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timeTaken = @(
    "2021-10-05 13:41:03.648234",
    "2021-10-05 13:41:03.648247",
    "2021-10-05 13:41:03.648251",
    "2021-10-05 13:41:03.648254",
    "2021-10-05 13:41:03.648258",
    "2021-10-05 13:41:03.648261",
    "2021-10-05 13:41:03.648264",
    "2021-10-05 13:41:03.648267",
    "2021-10-05 13:41:03.648271",
    "2021-10-05 13:41:03.648274",
    "2021-10-05 13:41:03.648277",
    "2021-10-05 13:41:03.648280",
    "2021-10-05 13:41:03.648283",
    "2021-10-05 13:41:03.648286",
    "2021-10-05 13:41:03.648289",
    "2021-10-05 13:41:03.648292",
    "2021-10-05 13:41:03.648295",
    "2021-10-05 13:41:03.648298",
    "2021-10-05 13:41:03.648301",
    "2021-10-05 13:41:03.648304",
    "2021-10-05 13:41:03.648307",
    "2021-10-05 13:41:03.648310",
    "2021-10-05 13:41:03.648313",
    "2021-10-05 13:41:03.648316",
    "2021-10-05 13:41:03.648319",
    "2021-10-05 13:41:03.648323",
    "2021-10-05 13:41:03.648326",
    "2021-10-05 13:41:03.648329",
    "2021-10-05 13:41:03.648332",
    "2021-10-05 13:41:03.648335",
    "2021-10-05 13:41:03.648338",
    "2021-10-05 13:41:03.648341",
    "2021-10-05 13:41:03.648344",
    "2021-10-05 13:41:03.648348",
    "2021-10-05 13:41:03.648351",
    "2021-10-05 13:41:03.648354",
    "2021-10-05 13:41:03.648357",
    "2021-10-05 13:41:03.648360",
    "2021-10-05 13:41:03.648363",
    "2021-10-05 13:41:03.648366",
    "2021-10-05 13:41:03.648369",
    "2021-10-05 13:41:03.648372",
    "2021-10-05 13:41:03.648375",
    "2021-10-05 13:41:03.648378",
    "2021-10-05 13:41:03.648381",
    "2021-10-05 13:41:03.648384",
    "2021-10-05 13:41:03.648387",
    "2021-10-05 13:41:03.648390",
    "2021-10-05 13:41:03.648393",
    "2021-10-05 13:41:03.648396",
    "2021-10-05 13:41:03.648399",
    "2021-10-05 13:41:03.648402",
    "2021-10-05 13:41:03.648405",
    "2021-10-05 13:41:03.648408",
    "2021-10-05 13:41:03.648411",
    "2021-10-05 13:41:03.648414",
    "2021-10-05 13:41:03.648417",
    "2021-10-05 13:41:03.648420",
    "2021-10-05 13:41:03.648423",
    "2021-10-05 13:41:03.648426",
    "2021-10-05 13:41:03.648429",
    "2021-10-05 13:41:03.648432",
    "2021-10-05 13:41:03.648435",
    "2021-10-05 13:41:03.648438",
    "2021-10-05 13:41:03.648443",
    "2021-10-05 13:41:03.648446",
    "2021-10-05 13:41:03.648449",
    "2021-10-05 13:41:03.648452",
    "2021-10-05 13:41:03.648455",
    "2021-10-05 13:41:03.648458",
    "2021-10-05 13:41:03.648461",
    "2021-10-05 13:41:03.648464",
    "2021-10-05 13:41:03.648467",
    "2021-10-05 13:41:03.648470",
    "2021-10-05 13:41:03.648473",
    "2021-10-05 13:41:03.648476",
    "2021-10-05 13:41:03.648481"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
